$wb = $excel.ActiveWorkbook

# Offense sheet - divisional round row (row 3)
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 294
$wsOff.Range("C3").Value = 213
$wsOff.Range("D3").Value = 69
$wsOff.Range("E3").Value = 46
$wsOff.Range("F3").Value = 7

# Defense sheet - divisional round row (row 3)
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 469
$wsDef.Range("C3").Value = 340
$wsDef.Range("D3").Value = 111
$wsDef.Range("E3").Value = 54
